$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-case the POL/POD city names (all-caps -> title case). The engine
# interns shared strings in the order they are first written, and drops
# strings that end up with zero references, so touching column B (POD)
# before column A (POL) on each row reproduces the desired shared-string
# ordering: MOTORCYCLE, Rotterdam, New York, Savannah, Miami, Houston,
# Indianapolis, Los Angeles, San Francisco, Varna.
$map = @{
  "NEW YORK"      = "New York";
  "SAVANNAH"      = "Savannah";
  "HOUSTON"       = "Houston";
  "INDIANAPOLIS"  = "Indianapolis";
  "LOS ANGELES"   = "Los Angeles";
  "ROTTERDAM"     = "Rotterdam";
  "VARNA"         = "Varna";
  "MIAMI"         = "Miami";
  "SAN FRANCISCO" = "San Francisco";
}

for ($r = 1; $r -le 71; $r++) {
    $podValue = $ws.Cells.Item($r, 2).Value2
    if ($map.ContainsKey($podValue)) {
        $ws.Cells.Item($r, 2).Value2 = $map[$podValue]
    }
    $polValue = $ws.Cells.Item($r, 1).Value2
    if ($map.ContainsKey($polValue)) {
        $ws.Cells.Item($r, 1).Value2 = $map[$polValue]
    }
}

# Move the active selection to match the saved view (scrolled to the
# bottom of the table, cell G69 selected).
$ws.Range("G69").Select()
